$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5 ("Entity Relationships"): fix punctuation in bullet 3 text and
# re-crop / reposition the diagram picture.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5Content = $s5.Shapes.Item(2)
$s5Content.TextFrame.TextRange.Paragraphs().Item(3).Text = "Firstly, One for User information like Username, Password, Email, Etc."

$s5Pic = $s5.Shapes.Item(3)
$s5Pic.PictureFormat.CropTop = 199.09125
$s5Pic.PictureFormat.CropRight = 375.0615
$s5Pic.Left = 471.7299837598425
$s5Pic.Top = 185.1707874015748
$s5Pic.Width = 384.1884876968504
$s5Pic.Height = 295.6583464566929

# ---------------------------------------------------------------------------
# Slide 6 ("Demonstration"): rewrite the body bullet list with the extra
# demonstration walk-through text, and grow the placeholder to fit.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6Content = $s6.Shapes.Item(2)
$s6Content.Height = 312.1187401574803
$s6Content.TextFrame.TextRange.Text = "Now onto a Demonstration…`rThat’s if we could easily but not so.`rSo, we’ll describe what “would” happen.`rWe would start by starting up / going to the site and login in or signing up`rThen from there showing the multiple pages and then final go to the NOTE page`rWhere we show how the note would work, like creating a new note, editing, saving, and deletion.`r"

# ---------------------------------------------------------------------------
# Restore the three deleted slides: insert 3 blank "Title and Content"
# slides right before the closing "Questions?" slide (which slides down to
# become the new last slide).
# ---------------------------------------------------------------------------
$questionsIndex = $p.Slides.Count
$p.Slides.Add($questionsIndex, 2) | Out-Null
$p.Slides.Add($questionsIndex, 2) | Out-Null
$p.Slides.Add($questionsIndex, 2) | Out-Null
